$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): rename headers, add new "note" header in C1, remove bold ---
$ws.Range("A1").Value = "codice_1_livello"
$ws.Range("B1").Value = "label_1_livello_it"
$ws.Range("C1").Value = "note"
$ws.Range("A1:C1").Font.Bold = $false

# --- New "note" column values marking deprecated/invalid items ---
$ws.Range("C15").Value = "Non Valido"
$ws.Range("C16").Value = "Non Valido"
$ws.Range("C20").Value = "Non Valido"
$ws.Range("C27").Value = "Non Valido"
$ws.Range("C29").Value = "Non Valido"
$ws.Range("C30").Value = "Non Valido"
$ws.Range("C31").Value = "Non Valido"

# --- Row 2: rename label, remove bold from B2 ---
$ws.Range("B2").Value = "Intestatario della Scheda"
$ws.Range("B2").Font.Bold = $false

# --- Text fixes (capitalization) ---
$ws.Range("B12").Value = "Zio / Zia (collaterale)"
$ws.Range("B21").Value = "Nipote (affine)"
$ws.Range("B22").Value = "Zio / Zia (affine)"

# --- Sheet view updates ---
$ws.Range("C32").Select()
$excel.ActiveWindow.ScrollRow = 20

$wb.Save()
